# Server done, test app add.
#
# - D2/D3 held the shared string "test"; it becomes the descriptive
#   "Корпус (test), этаж (test)" (both cells keep sharing the same
#   underlying string, so updating either Range.Value updates the shared
#   string text in place and both cells keep pointing at it).
# - Column D is widened to fit the new, longer text.
# - The last user selection on the sheet moves to E5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Корпус (test), этаж (test)"
$ws.Range("D3").Value = "Корпус (test), этаж (test)"

# Widen column D to show the longer text (closest reachable value to the
# 21.33203125-character stored width Excel's own autofit produced).
$ws.Columns.Item(4).ColumnWidth = 20.43

# Last active selection ends up on E5.
$ws.Range("E5").Select()
